# Add "Categories" and "Products" sheets (with sample data) to the workbook,
# appended after the existing "Purchases" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Categories sheet
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsCategories = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsCategories.Name = "Categories"

$categoriesHeaders = @("id", "name", "description")
for ($col = 1; $col -le $categoriesHeaders.Count; $col++) {
    $wsCategories.Cells.Item(1, $col).Value = $categoriesHeaders[$col - 1]
}

$categoriesRows = @(
    @("CAT001", "Decorative Laminates", "High-quality decorative laminates for furniture and interiors"),
    @("CAT002", "Industrial Laminates", "Durable laminates for industrial applications"),
    @("CAT003", "Compact Laminates", "High-pressure compact laminates for heavy-duty use"),
    @("CAT004", "Wood Finish Laminates", "Natural wood-finish laminate sheets"),
    @("CAT005", "Solid Color Laminates", "Solid color laminate sheets for modern designs")
)

$rowIndex = 2
foreach ($row in $categoriesRows) {
    for ($col = 1; $col -le $row.Count; $col++) {
        $wsCategories.Cells.Item($rowIndex, $col).Value = $row[$col - 1]
    }
    $rowIndex++
}

# ---------------------------------------------------------------------------
# Products sheet
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsProducts = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsProducts.Name = "Products"

$productsHeaders = @("id", "name", "categoryId", "categoryName", "length", "width", "thickness", "area", "unitCost", "colour", "weight", "notes")
for ($col = 1; $col -le $productsHeaders.Count; $col++) {
    $wsProducts.Cells.Item(1, $col).Value = $productsHeaders[$col - 1]
}

# Numeric columns (length, width, thickness, area, unitCost, weight) are real
# numbers; the rest stay text.
$productsRows = @(
    @("PROD001", "Premium Oak Laminate", "CAT004", "Wood Finish Laminates", 2440, 1220, 1, 2976800, 0.0025, "Oak Brown", 15000, "Popular wood finish, suitable for furniture"),
    @("PROD002", "Classic Walnut Laminate", "CAT004", "Wood Finish Laminates", 2440, 1220, 0.8, 2976800, 0.0028, "Dark Walnut", 14000, "Premium walnut finish"),
    @("PROD003", "Glossy White Laminate", "CAT005", "Solid Color Laminates", 2440, 1220, 1, 2976800, 0.002, "White", 15000, "High-gloss white finish"),
    @("PROD004", "Matte Black Laminate", "CAT005", "Solid Color Laminates", 2440, 1220, 1, 2976800, 0.0022, "Black", 15200, "Matte black finish for modern designs"),
    @("PROD005", "Industrial Grey Compact", "CAT003", "Compact Laminates", 3050, 1300, 12, 3965000, 0.0045, "Grey", 48000, "High-pressure compact laminate for heavy-duty applications"),
    @("PROD006", "Marble Effect Laminate", "CAT001", "Decorative Laminates", 2440, 1220, 1, 2976800, 0.0032, "White Marble", 15000, "Realistic marble pattern for premium interiors")
)

$rowIndex = 2
foreach ($row in $productsRows) {
    for ($col = 1; $col -le $row.Count; $col++) {
        $wsProducts.Cells.Item($rowIndex, $col).Value = $row[$col - 1]
    }
    $rowIndex++
}
